$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-01-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-31 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("741÷2=370, 1", $true, $false, $false, $false, $false, $true, 1, $false, "306÷4=76, 2", 2) | Out-Null
$d.Content.Find.Execute("435÷3=145, 0", $true, $false, $false, $false, $false, $true, 1, $false, "743÷6=123, 5", 2) | Out-Null
$d.Content.Find.Execute("894÷2=447, 0", $true, $false, $false, $false, $false, $true, 1, $false, "360÷4=90, 0", 2) | Out-Null
$d.Content.Find.Execute("764÷5=152, 4", $true, $false, $false, $false, $false, $true, 1, $false, "923÷2=461, 1", 2) | Out-Null
$d.Content.Find.Execute("928÷8=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "473÷2=236, 1", 2) | Out-Null
$d.Content.Find.Execute("208÷3=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "984÷4=246, 0", 2) | Out-Null
$d.Content.Find.Execute("373÷3=124, 1", $true, $false, $false, $false, $false, $true, 1, $false, "278÷9=30, 8", 2) | Out-Null
$d.Content.Find.Execute("140÷2=70, 0", $true, $false, $false, $false, $false, $true, 1, $false, "230÷3=76, 2", 2) | Out-Null
$d.Content.Find.Execute("363÷7=51, 6", $true, $false, $false, $false, $false, $true, 1, $false, "429÷9=47, 6", 2) | Out-Null
$d.Content.Find.Execute("275÷2=137, 1", $true, $false, $false, $false, $false, $true, 1, $false, "604÷2=302, 0", 2) | Out-Null
$d.Content.Find.Execute("713÷8=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "335÷2=167, 1", 2) | Out-Null
$d.Content.Find.Execute("558÷3=186, 0", $true, $false, $false, $false, $false, $true, 1, $false, "728÷3=242, 2", 2) | Out-Null
$d.Content.Find.Execute("921÷9=102, 3", $true, $false, $false, $false, $false, $true, 1, $false, "197÷7=28, 1", 2) | Out-Null
$d.Content.Find.Execute("207÷8=25, 7", $true, $false, $false, $false, $false, $true, 1, $false, "562÷9=62, 4", 2) | Out-Null
$d.Content.Find.Execute("751÷6=125, 1", $true, $false, $false, $false, $false, $true, 1, $false, "124÷4=31, 0", 2) | Out-Null
$d.Content.Find.Execute("921÷3=307, 0", $true, $false, $false, $false, $false, $true, 1, $false, "151÷4=37, 3", 2) | Out-Null
$d.Content.Find.Execute("465÷7=66, 3", $true, $false, $false, $false, $false, $true, 1, $false, "301÷9=33, 4", 2) | Out-Null
$d.Content.Find.Execute("408÷4=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "576÷3=192, 0", 2) | Out-Null
$d.Content.Find.Execute("554÷7=79, 1", $true, $false, $false, $false, $false, $true, 1, $false, "817÷5=163, 2", 2) | Out-Null
$d.Content.Find.Execute("711÷7=101, 4", $true, $false, $false, $false, $false, $true, 1, $false, "857÷6=142, 5", 2) | Out-Null
$d.Content.Find.Execute("681÷4=170, 1", $true, $false, $false, $false, $false, $true, 1, $false, "948÷2=474, 0", 2) | Out-Null
$d.Content.Find.Execute("782÷9=86, 8", $true, $false, $false, $false, $false, $true, 1, $false, "676÷4=169, 0", 2) | Out-Null
$d.Content.Find.Execute("244÷9=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "556÷4=139, 0", 2) | Out-Null
$d.Content.Find.Execute("100÷8=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "458÷6=76, 2", 2) | Out-Null
$d.Content.Find.Execute("529÷9=58, 7", $true, $false, $false, $false, $false, $true, 1, $false, "593÷3=197, 2", 2) | Out-Null
